$wb = $excel.ActiveWorkbook

# The "optimization_parameters" sheet had a stray row (labeled "Sheet",
# holding the values 3 and 4) that is no longer needed -- remove it
# entirely so the rows below shift up.
$paramsSheet = $wb.Worksheets.Item("optimization_parameters")
$paramsSheet.Rows.Item(16).Select() | Out-Null
$paramsSheet.Rows.Item(16).Delete() | Out-Null

# Wrap up the audit on the "optimization_diagnostics" sheet, which becomes
# the active tab/selection for the workbook.
$diagSheet = $wb.Worksheets.Item("optimization_diagnostics")
$diagSheet.Select() | Out-Null
